$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.369377255439758
$ws.Range("B1").Value = 4.013579368591309
$ws.Range("C1").Value = 3.341412305831909
$ws.Range("D1").Value = 2.263534545898438
$ws.Range("E1").Value = 0.8341880440711975
